$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Ccl11"
$ws.Cells.Item(2, 3).Value = "Ccr3"
$ws.Cells.Item(2, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 3.872739
$ws.Cells.Item(2, 8).Value = 7.745478
$ws.Cells.Item(2, 9).Value = 0.03299520440786341
$ws.Cells.Item(2, 10).Value = 0.02275854185403964
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.1626153333333333
$ws.Cells.Item(2, 14).Value = 0.487846
$ws.Cells.Item(2, 15).Value = 0.1293260700537641
$ws.Cells.Item(2, 16).Value = 0.1293260700537641
$ws.Cells.Item(2, 17).Value = 0.629766743398
$ws.Cells.Item(2, 18).Value = 3.778600460388
$ws.Cells.Item(2, 19).Value = 0.004267140116689611
$ws.Cells.Item(2, 20).Value = 0.002943272778137054

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Ccl11"
$ws.Cells.Item(3, 3).Value = "Ccr3"
$ws.Cells.Item(3, 4).Value = "Neutrophils"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 3.872739
$ws.Cells.Item(3, 8).Value = 7.745478
$ws.Cells.Item(3, 9).Value = 0.03299520440786341
$ws.Cells.Item(3, 10).Value = 0.02275854185403964
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.8767803333333334
$ws.Cells.Item(3, 14).Value = 2.630341
$ws.Cells.Item(3, 15).Value = 0.6972931302732585
$ws.Cells.Item(3, 16).Value = 0.6972931302732585
$ws.Cells.Item(3, 17).Value = 3.395541391333
$ws.Cells.Item(3, 18).Value = 20.373248347998
$ws.Cells.Item(3, 19).Value = 0.02300732936556509
$ws.Cells.Item(3, 20).Value = 0.01586937488985827

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Ccl11"
$ws.Cells.Item(4, 3).Value = "Ccr3"
$ws.Cells.Item(4, 4).Value = "Resolving-Mac"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 3.872739
$ws.Cells.Item(4, 8).Value = 7.745478
$ws.Cells.Item(4, 9).Value = 0.03299520440786341
$ws.Cells.Item(4, 10).Value = 0.02275854185403964
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.21801
$ws.Cells.Item(4, 14).Value = 0.65403
$ws.Cells.Item(4, 15).Value = 0.1733807996729775
$ws.Cells.Item(4, 16).Value = 0.1733807996729775
$ws.Cells.Item(4, 17).Value = 0.84429582939
$ws.Cells.Item(4, 18).Value = 5.06577497634
$ws.Cells.Item(4, 19).Value = 0.005720734925608708
$ws.Cells.Item(4, 20).Value = 0.00394589418604432

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Ccl11"
$ws.Cells.Item(5, 3).Value = "Ccr3"
$ws.Cells.Item(5, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 101.337382
$ws.Cells.Item(5, 8).Value = 304.012146
$ws.Cells.Item(5, 9).Value = 0.8633805772213771
$ws.Cells.Item(5, 10).Value = 0.8932790395734661
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.1626153333333333
$ws.Cells.Item(5, 14).Value = 0.487846
$ws.Cells.Item(5, 15).Value = 0.1293260700537641
$ws.Cells.Item(5, 16).Value = 0.1293260700537641
$ws.Cells.Item(5, 17).Value = 16.47901215305733
$ws.Cells.Item(5, 18).Value = 148.311109377516
$ws.Cells.Item(5, 19).Value = 0.1116576170127911
$ws.Cells.Item(5, 20).Value = 0.1155242676494372

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Ccl11"
$ws.Cells.Item(6, 3).Value = "Ccr3"
$ws.Cells.Item(6, 4).Value = "Neutrophils"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 101.337382
$ws.Cells.Item(6, 8).Value = 304.012146
$ws.Cells.Item(6, 9).Value = 0.8633805772213771
$ws.Cells.Item(6, 10).Value = 0.8932790395734661
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.8767803333333334
$ws.Cells.Item(6, 14).Value = 2.630341
$ws.Cells.Item(6, 15).Value = 0.6972931302732585
$ws.Cells.Item(6, 16).Value = 0.6972931302732585
$ws.Cells.Item(6, 17).Value = 88.85062356908735
$ws.Cells.Item(6, 18).Value = 799.6556121217861
$ws.Cells.Item(6, 19).Value = 0.6020293453078268
$ws.Cells.Item(6, 20).Value = 0.6228773377116721

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Ccl11"
$ws.Cells.Item(7, 3).Value = "Ccr3"
$ws.Cells.Item(7, 4).Value = "Resolving-Mac"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 101.337382
$ws.Cells.Item(7, 8).Value = 304.012146
$ws.Cells.Item(7, 9).Value = 0.8633805772213771
$ws.Cells.Item(7, 10).Value = 0.8932790395734661
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.21801
$ws.Cells.Item(7, 14).Value = 0.65403
$ws.Cells.Item(7, 15).Value = 0.1733807996729775
$ws.Cells.Item(7, 16).Value = 0.1733807996729775
$ws.Cells.Item(7, 17).Value = 22.09256264982
$ws.Cells.Item(7, 18).Value = 198.83306384838
$ws.Cells.Item(7, 19).Value = 0.1496936149007592
$ws.Cells.Item(7, 20).Value = 0.1548774342123568

# Row 8
$ws.Cells.Item(8, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(8, 2).Value = "Ccl11"
$ws.Cells.Item(8, 3).Value = "Ccr3"
$ws.Cells.Item(8, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 2.422325
$ws.Cells.Item(8, 8).Value = 7.266975
$ws.Cells.Item(8, 9).Value = 0.02063787632403778
$ws.Cells.Item(8, 10).Value = 0.02135255625150052
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.1626153333333333
$ws.Cells.Item(8, 14).Value = 0.487846
$ws.Cells.Item(8, 15).Value = 0.1293260700537641
$ws.Cells.Item(8, 16).Value = 0.1293260700537641
$ws.Cells.Item(8, 17).Value = 0.3939071873166667
$ws.Cells.Item(8, 18).Value = 3.54516468585
$ws.Cells.Item(8, 19).Value = 0.00266901543924343
$ws.Cells.Item(8, 20).Value = 0.002761442185608495

# Row 9
$ws.Cells.Item(9, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(9, 2).Value = "Ccl11"
$ws.Cells.Item(9, 3).Value = "Ccr3"
$ws.Cells.Item(9, 4).Value = "Neutrophils"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 2.422325
$ws.Cells.Item(9, 8).Value = 7.266975
$ws.Cells.Item(9, 9).Value = 0.02063787632403778
$ws.Cells.Item(9, 10).Value = 0.02135255625150052
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.8767803333333334
$ws.Cells.Item(9, 14).Value = 2.630341
$ws.Cells.Item(9, 15).Value = 0.6972931302732585
$ws.Cells.Item(9, 16).Value = 0.6972931302732585
$ws.Cells.Item(9, 17).Value = 2.123846920941667
$ws.Cells.Item(9, 18).Value = 19.114622288475
$ws.Cells.Item(9, 19).Value = 0.01439064938418067
$ws.Cells.Item(9, 20).Value = 0.01488899078794463

# Row 10
$ws.Cells.Item(10, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(10, 2).Value = "Ccl11"
$ws.Cells.Item(10, 3).Value = "Ccr3"
$ws.Cells.Item(10, 4).Value = "Resolving-Mac"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 2.422325
$ws.Cells.Item(10, 8).Value = 7.266975
$ws.Cells.Item(10, 9).Value = 0.02063787632403778
$ws.Cells.Item(10, 10).Value = 0.02135255625150052
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.21801
$ws.Cells.Item(10, 14).Value = 0.65403
$ws.Cells.Item(10, 15).Value = 0.1733807996729775
$ws.Cells.Item(10, 16).Value = 0.1733807996729775
$ws.Cells.Item(10, 17).Value = 0.52809107325
$ws.Cells.Item(10, 18).Value = 4.75281965925
$ws.Cells.Item(10, 19).Value = 0.003578211500613679
$ws.Cells.Item(10, 20).Value = 0.003702123277947393

# Row 11
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Ccl11"
$ws.Cells.Item(11, 3).Value = "Ccr3"
$ws.Cells.Item(11, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 7.912825
$ws.Cells.Item(11, 8).Value = 15.82565
$ws.Cells.Item(11, 9).Value = 0.06741618227271494
$ws.Cells.Item(11, 10).Value = 0.04650051525449849
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.1626153333333333
$ws.Cells.Item(11, 14).Value = 0.487846
$ws.Cells.Item(11, 15).Value = 0.1293260700537641
$ws.Cells.Item(11, 16).Value = 0.1293260700537641
$ws.Cells.Item(11, 17).Value = 1.286746674983333
$ws.Cells.Item(11, 18).Value = 7.7204800499
$ws.Cells.Item(11, 19).Value = 0.008718669911358464
$ws.Cells.Item(11, 20).Value = 0.0060137288933394

# Row 12
$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Ccl11"
$ws.Cells.Item(12, 3).Value = "Ccr3"
$ws.Cells.Item(12, 4).Value = "Neutrophils"
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 7.912825
$ws.Cells.Item(12, 8).Value = 15.82565
$ws.Cells.Item(12, 9).Value = 0.06741618227271494
$ws.Cells.Item(12, 10).Value = 0.04650051525449849
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 0.8767803333333334
$ws.Cells.Item(12, 14).Value = 2.630341
$ws.Cells.Item(12, 15).Value = 0.6972931302732585
$ws.Cells.Item(12, 16).Value = 0.6972931302732585
$ws.Cells.Item(12, 17).Value = 6.937809341108333
$ws.Cells.Item(12, 18).Value = 41.62685604665
$ws.Cells.Item(12, 19).Value = 0.04700884076801395
$ws.Cells.Item(12, 20).Value = 0.03242448984112866

# Row 13
$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Ccl11"
$ws.Cells.Item(13, 3).Value = "Ccr3"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 7.912825
$ws.Cells.Item(13, 8).Value = 15.82565
$ws.Cells.Item(13, 9).Value = 0.06741618227271494
$ws.Cells.Item(13, 10).Value = 0.04650051525449849
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.21801
$ws.Cells.Item(13, 14).Value = 0.65403
$ws.Cells.Item(13, 15).Value = 0.1733807996729775
$ws.Cells.Item(13, 16).Value = 0.1733807996729775
$ws.Cells.Item(13, 17).Value = 1.72507497825
$ws.Cells.Item(13, 18).Value = 10.3504498695
$ws.Cells.Item(13, 19).Value = 0.01168867159334252
$ws.Cells.Item(13, 20).Value = 0.008062296520030435

# Row 14
$ws.Cells.Item(14, 1).Value = "Neutrophils"
$ws.Cells.Item(14, 2).Value = "Ccl11"
$ws.Cells.Item(14, 3).Value = "Ccr3"
$ws.Cells.Item(14, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 0.6666666666666666
$ws.Cells.Item(14, 7).Value = 0.5045936666666666
$ws.Cells.Item(14, 8).Value = 1.513781
$ws.Cells.Item(14, 9).Value = 0.004299068740387607
$ws.Cells.Item(14, 10).Value = 0.004447943464089625
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 12).Value = 0.6666666666666666
$ws.Cells.Item(14, 13).Value = 0.1626153333333333
$ws.Cells.Item(14, 14).Value = 0.487846
$ws.Cells.Item(14, 15).Value = 0.1293260700537641
$ws.Cells.Item(14, 16).Value = 0.1293260700537641
$ws.Cells.Item(14, 17).Value = 0.08205466730288888
$ws.Cells.Item(14, 18).Value = 0.7384920057259999
$ws.Cells.Item(14, 19).Value = 0.0005559816650853153
$ws.Cells.Item(14, 20).Value = 0.0005752350480320371

# Row 15
$ws.Cells.Item(15, 1).Value = "Neutrophils"
$ws.Cells.Item(15, 2).Value = "Ccl11"
$ws.Cells.Item(15, 3).Value = "Ccr3"
$ws.Cells.Item(15, 4).Value = "Neutrophils"
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 0.6666666666666666
$ws.Cells.Item(15, 7).Value = 0.5045936666666666
$ws.Cells.Item(15, 8).Value = 1.513781
$ws.Cells.Item(15, 9).Value = 0.004299068740387607
$ws.Cells.Item(15, 10).Value = 0.004447943464089625
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 0.8767803333333334
$ws.Cells.Item(15, 14).Value = 2.630341
$ws.Cells.Item(15, 15).Value = 0.6972931302732585
$ws.Cells.Item(15, 16).Value = 0.6972931302732585
$ws.Cells.Item(15, 17).Value = 0.4424178032578889
$ws.Cells.Item(15, 18).Value = 3.981760229320999
$ws.Cells.Item(15, 19).Value = 0.002997711099244789
$ws.Cells.Item(15, 20).Value = 0.003101520421353535

# Row 16
$ws.Cells.Item(16, 1).Value = "Neutrophils"
$ws.Cells.Item(16, 2).Value = "Ccl11"
$ws.Cells.Item(16, 3).Value = "Ccr3"
$ws.Cells.Item(16, 4).Value = "Resolving-Mac"
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 0.6666666666666666
$ws.Cells.Item(16, 7).Value = 0.5045936666666666
$ws.Cells.Item(16, 8).Value = 1.513781
$ws.Cells.Item(16, 9).Value = 0.004299068740387607
$ws.Cells.Item(16, 10).Value = 0.004447943464089625
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 0.21801
$ws.Cells.Item(16, 14).Value = 0.65403
$ws.Cells.Item(16, 15).Value = 0.1733807996729775
$ws.Cells.Item(16, 16).Value = 0.1733807996729775
$ws.Cells.Item(16, 17).Value = 0.11000646527
$ws.Cells.Item(16, 18).Value = 0.9900581874299998
$ws.Cells.Item(16, 19).Value = 0.0007453759760575032
$ws.Cells.Item(16, 20).Value = 0.0007711879947040526

# Row 17
$ws.Cells.Item(17, 1).Value = "Resolving-Mac"
$ws.Cells.Item(17, 2).Value = "Ccl11"
$ws.Cells.Item(17, 3).Value = "Ccr3"
$ws.Cells.Item(17, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 1.322919333333333
$ws.Cells.Item(17, 8).Value = 3.968758
$ws.Cells.Item(17, 9).Value = 0.01127109103361929
$ws.Cells.Item(17, 10).Value = 0.01166140360240577
$ws.Cells.Item(17, 11).Value = 2
$ws.Cells.Item(17, 12).Value = 0.6666666666666666
$ws.Cells.Item(17, 13).Value = 0.1626153333333333
$ws.Cells.Item(17, 14).Value = 0.487846
$ws.Cells.Item(17, 15).Value = 0.1293260700537641
$ws.Cells.Item(17, 16).Value = 0.1293260700537641
$ws.Cells.Item(17, 17).Value = 0.2151269683631111
$ws.Cells.Item(17, 18).Value = 1.936142715268
$ws.Cells.Item(17, 19).Value = 0.001457645908596201
$ws.Cells.Item(17, 20).Value = 0.001508123499209947

# Row 18
$ws.Cells.Item(18, 1).Value = "Resolving-Mac"
$ws.Cells.Item(18, 2).Value = "Ccl11"
$ws.Cells.Item(18, 3).Value = "Ccr3"
$ws.Cells.Item(18, 4).Value = "Neutrophils"
$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 1.322919333333333
$ws.Cells.Item(18, 8).Value = 3.968758
$ws.Cells.Item(18, 9).Value = 0.01127109103361929
$ws.Cells.Item(18, 10).Value = 0.01166140360240577
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 12).Value = 1
$ws.Cells.Item(18, 13).Value = 0.8767803333333334
$ws.Cells.Item(18, 14).Value = 2.630341
$ws.Cells.Item(18, 15).Value = 0.6972931302732585
$ws.Cells.Item(18, 16).Value = 0.6972931302732585
$ws.Cells.Item(18, 17).Value = 1.159909654053111
$ws.Cells.Item(18, 18).Value = 10.439186886478
$ws.Cells.Item(18, 19).Value = 0.00785925434842725
$ws.Cells.Item(18, 20).Value = 0.008131416621301374

# Row 19
$ws.Cells.Item(19, 1).Value = "Resolving-Mac"
$ws.Cells.Item(19, 2).Value = "Ccl11"
$ws.Cells.Item(19, 3).Value = "Ccr3"
$ws.Cells.Item(19, 4).Value = "Resolving-Mac"
$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 1.322919333333333
$ws.Cells.Item(19, 8).Value = 3.968758
$ws.Cells.Item(19, 9).Value = 0.01127109103361929
$ws.Cells.Item(19, 10).Value = 0.01166140360240577
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 0.21801
$ws.Cells.Item(19, 14).Value = 0.65403
$ws.Cells.Item(19, 15).Value = 0.1733807996729775
$ws.Cells.Item(19, 16).Value = 0.1733807996729775
$ws.Cells.Item(19, 17).Value = 0.28840964386
$ws.Cells.Item(19, 18).Value = 2.59568679474
$ws.Cells.Item(19, 19).Value = 0.001954190776595838
$ws.Cells.Item(19, 20).Value = 0.002021863481894453

Write-Host "done"